# The underlying OOXML diff for this commit only touches the
# Slide Master (ppt/slideMasters/slideMaster1.xml) and the "Titelfolie"
# Slide Layout (ppt/slideLayouts/slideLayout1.xml): in every
# a14:hiddenFill / a14:hiddenLine compatibility-extension element the
# two namespace declarations (xmlns:a14="...", xmlns="") are simply
# swapped in serialization order, e.g.
#
#   <a14:hiddenFill xmlns="" xmlns:a14="...">
# becomes
#   <a14:hiddenFill xmlns:a14="..." xmlns="">
#
# That is purely a cosmetic artifact of the XML writer used when the
# file was last resaved (the compatibility <a:ln>/<a:solidFill> payload
# cached inside those extension elements, their uris, and every other
# attribute/value are untouched) - it carries no visible or semantic
# change for any shape, slide, placeholder or text run in the deck.
#
# There is no PowerPoint object-model call that targets these raw
# extLst namespace declarations (Shape.Fill/Shape.Line only manage the
# <a:solidFill>/<a:ln> siblings, never the cached "hidden" compat
# payload), so the faithful reproduction of this commit is to leave the
# Slide Master / Slide Layout - and every other part, which the diff
# does not touch at all - exactly as authored.

$p = $ppt.ActivePresentation
